# Apply "calculation of new indicators" changes across the measures/metrics
# sheets of the workbook: renumber indicator codes and recalculate values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# SCHEME_MEASURES: MQMS0x -> MQME00x (same descriptions / values)
# ---------------------------------------------------------------------------
$wsScheme = $wb.Worksheets.Item("SCHEME_MEASURES")
$wsScheme.Cells.Item(2, 1).Value = "MQME001"
$wsScheme.Cells.Item(3, 1).Value = "MQME002"
$wsScheme.Cells.Item(4, 1).Value = "MQME003"
$wsScheme.Cells.Item(5, 1).Value = "MQME004"
$wsScheme.Cells.Item(6, 1).Value = "MQME005"

# ---------------------------------------------------------------------------
# METADATA_ISSUES: renumber rule codes in column A only
# ---------------------------------------------------------------------------
$wsIssues = $wb.Worksheets.Item("METADATA_ISSUES")
foreach ($r in 2..3) {
    $wsIssues.Cells.Item($r, 1).Value = "MQME012"
}
foreach ($r in 4..9) {
    $wsIssues.Cells.Item($r, 1).Value = "MQME014"
}
foreach ($r in 10..39) {
    $wsIssues.Cells.Item($r, 1).Value = "MQME008"
}
foreach ($r in 40..43) {
    $wsIssues.Cells.Item($r, 1).Value = "MQME011"
}

# ---------------------------------------------------------------------------
# METADATA_MEASURES: drop the "Total number of columns" row, renumber and
# recompute the remaining two rows
# ---------------------------------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("METADATA_MEASURES")
$wsMeasures.Rows.Item(4).Delete()

$wsMeasures.Cells.Item(2, 1).Value = "MQME006"
$wsMeasures.Cells.Item(2, 2).Value = "Total number of length-required columns"
$wsMeasures.Cells.Item(2, 3).Value = 139

$wsMeasures.Cells.Item(3, 1).Value = "MQME007"
$wsMeasures.Cells.Item(3, 2).Value = "Total number of NUMBER columns"
$wsMeasures.Cells.Item(3, 3).Value = 399

# ---------------------------------------------------------------------------
# METADATA_METRICS: renumber IQME0x -> MQID0xx, reorder/relabel rows and add
# four new indicator rows
# ---------------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("METADATA_METRICS")

$metricsData = @(
    @("MQID001", "Table names in singular", "97.06%"),
    @("MQID002", "Table with recommended name length", "100.00%"),
    @("MQID003", "Columns with correct prefixes", "99.06%"),
    @("MQID004", "Columns with recommended name size", "100.00%"),
    @("MQID005", "Columns with comments", "95.31%"),
    @("MQID006", "Table with standard PK prefixes", "100.00%"),
    @("MQID007", "Table with standard FK prefixes", "100.00%"),
    @("MQID008", "Table with standard UK prefixes", "73.33%"),
    @("MQID009", "NUMBER columns with valid scale", "100.00%"),
    @("MQID010", "Columns with valid num_distinct", "100.00%"),
    @("MQID011", "Columns with valid num_nulls", "100.00%")
)

$row = 2
foreach ($entry in $metricsData) {
    $wsMetrics.Cells.Item($row, 1).Value = $entry[0]
    $wsMetrics.Cells.Item($row, 2).Value = $entry[1]
    $wsMetrics.Cells.Item($row, 3).Value = "'" + $entry[2]
    $row = $row + 1
}
